# Adapt column header formatting to respective input file names (#7)
#  - rename header suffixes "_old" -> "_FV2404" and "_new" -> "_FV2410"
#  - turn the data range into an Excel Table (ListObject)
#  - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2404 = "_FV2404"
$fv2410 = "_FV2410"

# Rename header cells in row 1 (A1:U1), replacing the "_old"/"_new" suffixes
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value2

    if ($text.EndsWith($oldSuffix)) {
        $base = $text.Substring(0, $text.Length - $oldSuffix.Length)
        $cell.Value = $base + $fv2404
    }
    elseif ($text.EndsWith($newSuffix)) {
        $base = $text.Substring(0, $text.Length - $newSuffix.Length)
        $cell.Value = $base + $fv2410
    }
}

# Convert the used range into a native Excel table ("Table1")
$tableRange = $ws.UsedRange
$listObject = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# Freeze the header row (split below row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
